$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Sheets("展览")
$ws1.Range("F5").Value2 = 15439
$ws1.Range("F9").Value2 = 15363
$ws1.Range("F10").Value2 = 49
$ws1.Range("F11").Value2 = 8958
$ws1.Range("F18").Value2 = 194
$ws1.Range("F26").Value2 = 13
$ws1.Range("F28").Value2 = 76
$ws1.Range("F30").Value2 = 38
$ws1.Range("F34").Value2 = 242
$ws1.Range("F35").Value2 = 302
$ws1.Range("F36").Value2 = 441
$ws1.Range("F38").Value2 = 5480

# Sheet "全部类型" (All types, combined)
$ws4 = $wb.Sheets("全部类型")
$ws4.Range("F5").Value2 = 15439
$ws4.Range("F9").Value2 = 15363
$ws4.Range("F10").Value2 = 49
$ws4.Range("F11").Value2 = 8958
$ws4.Range("F18").Value2 = 194
$ws4.Range("F26").Value2 = 13
$ws4.Range("F28").Value2 = 76
$ws4.Range("F30").Value2 = 38
$ws4.Range("F36").Value2 = 242
$ws4.Range("F37").Value2 = 302
$ws4.Range("F38").Value2 = 441
$ws4.Range("F40").Value2 = 5480
